$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 8).Value = [double]"0.09848601862732026"
$ws.Cells.Item(2, 9).Value = [double]"1.652238118382965"
$ws.Cells.Item(3, 8).Value = [double]"1.321290021644665e-41"
$ws.Cells.Item(3, 9).Value = [double]"-13.51240738922745"
$ws.Cells.Item(4, 8).Value = [double]"2.039558748889865e-77"
$ws.Cells.Item(4, 9).Value = [double]"18.62434938601493"
$ws.Cells.Item(5, 8).Value = [double]"2.317154994930742e-17"
$ws.Cells.Item(5, 9).Value = [double]"8.476677804302749"
$ws.Cells.Item(6, 8).Value = [double]"5.277542374521448e-25"
$ws.Cells.Item(6, 9).Value = [double]"10.32776590728914"
$ws.Cells.Item(7, 8).Value = [double]"0.0002792762503104448"
$ws.Cells.Item(7, 9).Value = [double]"3.633801040501866"
$ws.Cells.Item(8, 8).Value = [double]"2.216290846907834e-174"
$ws.Cells.Item(8, 9).Value = [double]"28.15277166673561"
$ws.Cells.Item(9, 8).Value = [double]"3.223078845225822e-97"
$ws.Cells.Item(9, 9).Value = [double]"-20.92421789370874"
$ws.Cells.Item(10, 6).Value = [double]"4.090125119497237e-14"
$ws.Cells.Item(10, 8).Value = [double]"7.779370280598121e-51"
$ws.Cells.Item(10, 9).Value = [double]"-14.99615821514027"
$ws.Cells.Item(11, 8).Value = [double]"2.148157955288719e-50"
$ws.Cells.Item(11, 9).Value = [double]"-14.92857245517089"
$ws.Cells.Item(12, 8).Value = [double]"1.405045832173766e-22"
$ws.Cells.Item(12, 9).Value = [double]"9.777598595581384"
$ws.Cells.Item(13, 8).Value = [double]"2.330551571583976e-290"
$ws.Cells.Item(13, 9).Value = [double]"36.41653484576065"
$ws.Cells.Item(14, 8).Value = [double]"8.225420312876559e-96"
$ws.Cells.Item(14, 9).Value = [double]"-20.7691777957254"
$ws.Cells.Item(15, 8).Value = [double]"8.551553098984214e-38"
$ws.Cells.Item(15, 9).Value = [double]"-12.8504462250916"
$ws.Cells.Item(16, 8).Value = [double]"3.683496002766283e-26"
$ws.Cells.Item(16, 9).Value = [double]"-10.58015468445819"
$ws.Cells.Item(17, 8).Value = [double]"1.308917393152356e-31"
$ws.Cells.Item(17, 9).Value = [double]"11.69774092658699"
$ws.Cells.Item(18, 8).Value = [double]"0.0009446089762037671"
$ws.Cells.Item(18, 9).Value = [double]"-3.306525069182445"
$ws.Cells.Item(19, 8).Value = [double]"2.45394121593645e-92"
$ws.Cells.Item(19, 9).Value = [double]"20.38123324251259"
$ws.Cells.Item(20, 8).Value = [double]"1.504752556286769e-101"
$ws.Cells.Item(20, 9).Value = [double]"21.39445581604381"
$ws.Cells.Item(21, 8).Value = [double]"1.04114863183972e-112"
$ws.Cells.Item(21, 9).Value = [double]"22.56125994674017"
$ws.Cells.Item(22, 6).Value = [double]"0.005424722780159436"
$ws.Cells.Item(22, 8).Value = [double]"0.0006815864733909085"
$ws.Cells.Item(22, 9).Value = [double]"3.39688081113558"
$ws.Cells.Item(23, 4).Value = [double]"4.133503021397445e-26"
$ws.Cells.Item(23, 6).Value = [double]"4.31305663103952e-26"
$ws.Cells.Item(23, 8).Value = [double]"1.281799887524342e-104"
$ws.Cells.Item(23, 9).Value = [double]"21.72161949362362"
$ws.Cells.Item(24, 8).Value = [double]"3.341975839169662e-54"
$ws.Cells.Item(24, 9).Value = [double]"-15.50239518950835"
$ws.Cells.Item(25, 8).Value = [double]"1.344667035640849e-19"
$ws.Cells.Item(25, 9).Value = [double]"-9.056691940797204"
$ws.Cells.Item(26, 8).Value = [double]"6.560262656697334e-28"
$ws.Cells.Item(26, 9).Value = [double]"-10.95115920964167"
$ws.Cells.Item(27, 8).Value = [double]"2.169412589807753e-20"
$ws.Cells.Item(27, 9).Value = [double]"9.253656569601779"
$ws.Cells.Item(28, 8).Value = [double]"1.071242209843778e-158"
$ws.Cells.Item(28, 9).Value = [double]"26.84110715351077"
$ws.Cells.Item(29, 8).Value = [double]"3.66317047257038e-32"
$ws.Cells.Item(29, 9).Value = [double]"-11.80533776269232"
$ws.Cells.Item(30, 8).Value = [double]"1.388731868567314e-12"
$ws.Cells.Item(30, 9).Value = [double]"-7.08517052951575"
$ws.Cells.Item(31, 8).Value = [double]"1.502438606793471e-11"
$ws.Cells.Item(31, 9).Value = [double]"-6.747664496003302"
$ws.Cells.Item(32, 8).Value = [double]"2.894783167268939e-27"
$ws.Cells.Item(32, 9).Value = [double]"-10.81588711661282"
$ws.Cells.Item(33, 6).Value = [double]"4.641747559277946e-07"
$ws.Cells.Item(33, 8).Value = [double]"5.176275505456987e-59"
$ws.Cells.Item(33, 9).Value = [double]"-16.19839225683151"
$ws.Cells.Item(34, 8).Value = [double]"3.838287071939599e-12"
$ws.Cells.Item(34, 9).Value = [double]"6.943011299358002"
$ws.Cells.Item(35, 8).Value = [double]"0.0176184553557072"
$ws.Cells.Item(35, 9).Value = [double]"2.37354061303342"
$ws.Cells.Item(36, 8).Value = [double]"0.0008382126776864868"
$ws.Cells.Item(36, 9).Value = [double]"3.339857898741211"
$ws.Cells.Item(37, 6).Value = [double]"0.2074807101666616"
$ws.Cells.Item(37, 8).Value = [double]"0.1899958221220994"
$ws.Cells.Item(37, 9).Value = [double]"1.310591471503704"
$ws.Cells.Item(38, 8).Value = [double]"4.176791926354502e-59"
$ws.Cells.Item(38, 9).Value = [double]"16.21158168338971"
$ws.Cells.Item(39, 8).Value = [double]"7.045463199321931e-43"
$ws.Cells.Item(39, 9).Value = [double]"-13.72650113197182"
$ws.Cells.Item(40, 8).Value = [double]"5.983863907945237e-21"
$ws.Cells.Item(40, 9).Value = [double]"-9.39028409075843"
$ws.Cells.Item(41, 8).Value = [double]"1.348493453021991e-22"
$ws.Cells.Item(41, 9).Value = [double]"-9.781756727778122"
$ws.Cells.Item(42, 8).Value = [double]"1.668480096421356e-05"
$ws.Cells.Item(42, 9).Value = [double]"4.30518265926056"
$ws.Cells.Item(43, 8).Value = [double]"1.125374071928214e-124"
$ws.Cells.Item(43, 9).Value = [double]"23.74899321397551"
$ws.Cells.Item(44, 8).Value = [double]"4.415798392557715e-82"
$ws.Cells.Item(44, 9).Value = [double]"-19.19082428022793"
$ws.Cells.Item(45, 8).Value = [double]"1.70807036917321e-28"
$ws.Cells.Item(45, 9).Value = [double]"-11.07237794302584"
$ws.Cells.Item(46, 8).Value = [double]"6.831209819502281e-17"
$ws.Cells.Item(46, 9).Value = [double]"-8.34991502017291"
